$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.442.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -3.22%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.839.37'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -3.58%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.71'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.94%  '

$ws.Range('E7').Value = '  -2.16%  '

$ws.Range('E8').Value = '  +0.20%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.746'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.36%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.174'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.88%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.07'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.14%  '

$ws.Range('E12').Value = '  +0.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.31'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.459.40'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.02'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.11%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.857.13'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.77'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.61%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.20'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -6.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.350.08'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '435.99'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.29%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.70'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '93.79'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.26'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.88%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.85'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.05%  '

$ws.Range('E26').Value = '  -9.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.16'
$ws.Range('D27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.94'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.23%  '

$ws.Range('E29').Value = '  -1.68%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.03'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.58%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.16'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.44'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.60%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '48.01'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.82%  '

$ws.Range('E34').Value = '  -4.43%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '69.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0973'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.57%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '628.72'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -7.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.422'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.92%  '

$ws.Range('E39').Value = '  -0.84%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.16%  '

$ws.Range('E41').Value = '  +0.10%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.28'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.93%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +24.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0467'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.00'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.69'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.12%  '

$ws.Range('E47').Value = '  -4.22%  '

$ws.Range('E48').Value = '  -15.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.824.79'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.91%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.24'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.49%  '

$ws.Range('E51').Value = '  +0.86%  '

